$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DW")
$ws.Activate()

# Add a new entry (row 16) to the table: "Remove Spaces from String"
$ws.Range("A16").Value = "Remove Spaces from string"
$ws.Range("B16").Value = "String"
$ws.Range("C16").Value = "Easy"
$ws.Range("D16").Value = "https://www.geeksforgeeks.org/remove-spaces-from-a-given-string/"
$ws.Range("E16").Value = "Traverse and only add no space characters in a list and then join the list to a empty string"

# Match formatting used by the rest of the table: "Neutral" cell style with a
# left/right border for A:C (copied from the row above so the exact same
# style/border combination is reused), and the bordered "Neutral" look for E.
# The Link column (D) is intentionally left unformatted, same as the source.
$ws.Range("A15:C15").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("E16").Style = "Neutral"
$ws.Range("E16").HorizontalAlignment = -4131       # xlLeft

$ws.Rows.Item(16).RowHeight = 16

$ws.Range("E16").Select() | Out-Null
